# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape. Values that look numeric but must stay text (e.g. "2.00",
# "0.450") are entered with a leading apostrophe to force text storage,
# then the cell Style is reset to "Normal" so no stray number format/
# quote-prefix style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '62.986.64'
$ws.Range("E2").Value = '  -0.37%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.067.52'
$ws.Range("E3").Value = '  -0.83%  '

# Row 5: BNB
$ws.Range("D5").Value = '''536.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.36%  '

# Row 6: Solana
$ws.Range("D6").Value = '''132.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.49%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.18%  '

# Row 8: LidoStakedEther
$ws.Range("D8").Value = '3.058.77'
$ws.Range("E8").Value = '  -0.81%  '

# Row 9: XRP
$ws.Range("D9").Value = '''0.492'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.86%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.153'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.71%  '

# Row 11: Toncoin
$ws.Range("D11").Value = '''6.11'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.22%  '

# Row 12: Cardano
$ws.Range("D12").Value = '''0.450'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.12%  '

# Row 13: ShibaInu
$ws.Range("E13").Value = '  +2.01%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''34.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.18%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.563.23'
$ws.Range("E15").Value = '  -0.38%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '63.007.67'
$ws.Range("E16").Value = '  -0.25%  '

# Row 17: TRON
$ws.Range("E17").Value = '  -0.33%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '3.070.33'
$ws.Range("E18").Value = '  -0.44%  '

# Row 19: Polkadot
$ws.Range("D19").Value = '''6.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.37%  '

# Row 20: BitcoinCash
$ws.Range("D20").Value = '''481.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.65%  '

# Row 21: Chainlink
$ws.Range("D21").Value = '''13.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.95%  '

# Row 22: Polygon
$ws.Range("D22").Value = '''0.692'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.32%  '

# Row 23: Uniswap
$ws.Range("D23").Value = '''7.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.18%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '''78.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.10%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range("D25").Value = '''12.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.92%  '

# Row 26: Dai
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.09%  '

# Row 27: PancakeSwap
$ws.Range("E27").Value = '  -2.64%  '

# Row 28: RenderToken
$ws.Range("D28").Value = '''8.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.06%  '

# Row 29: FirstDigitalUSD
$ws.Range("D29").Value = '''0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.08%  '

# Row 30: EthereumClassic
$ws.Range("D30").Value = '''25.90'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.57%  '

# Row 31: ImmutableX
$ws.Range("D31").Value = '''1.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.67%  '

# Row 32: Mantle
$ws.Range("E32").Value = '  -1.02%  '

# Row 33: Stacks
$ws.Range("D33").Value = '''2.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.90%  '

# Row 34: OKB
$ws.Range("D34").Value = '''56.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.69%  '

# Row 35: NEARProtocol
$ws.Range("D35").Value = '''5.32'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.91%  '

# Row 36: Filecoin
$ws.Range("D36").Value = '''5.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.01%  '

# Row 37: Bittensor
$ws.Range("D37").Value = '''479.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.62%  '

# Row 38: VeChain
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.095.69'
$ws.Range("E38").Value = '  +0.40%  '

# Row 39: Maker
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.0393'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.15%  '

# Row 40: Hedera
$ws.Range("D40").Value = '''0.0792'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.79%  '

# Row 41: Kaspa
$ws.Range("E41").Value = '  -2.69%  '

# Row 42: Cosmos
$ws.Range("D42").Value = '''8.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.93%  '

# Row 43: dogwifhat
$ws.Range("E43").Value = '  -2.20%  '

# Row 44: TheGraph
$ws.Range("D44").Value = '''0.251'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.83%  '

# Row 46: PEPE
$ws.Range("D46").Value = '0.0₃0536'
$ws.Range("E46").Value = '  +7.48%  '

# Row 47: Monero
$ws.Range("D47").Value = '''121.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.12%  '

# Row 48: Fetch.AI
$ws.Range("D48").Value = '''2.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.91%  '

# Row 49: InjectiveProtocol
$ws.Range("D49").Value = '''24.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.41%  '

# Row 50: Stellar
$ws.Range("E50").Value = '  +0.53%  '

# Row 51: ThetaToken
$ws.Range("D51").Value = '''2.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.20%  '
